$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: translate Spanish column headers to snake_case English names ---
$ws.Cells.Item(1, 1).Value = 'mx_state'
$ws.Cells.Item(1, 2).Value = 'mx_municipality'
$ws.Cells.Item(1, 3).Value = 'n_matriculas'
$ws.Cells.Item(1, 4).Value = 'pct_matriculas'

# --- Title-case the Spanish connector words (de/del/la/las/los/el/y) in state/municipality names ---
$ws.Cells.Item(8, 2).Value = 'Pabellón De Arteaga'
$ws.Cells.Item(9, 2).Value = 'Rincón De Romos'
$ws.Cells.Item(10, 2).Value = 'San Francisco De Los Romo'
$ws.Cells.Item(11, 2).Value = 'San José De Gracia'
$ws.Cells.Item(16, 2).Value = 'Playas De Rosarito'
$ws.Cells.Item(40, 2).Value = 'Amatenango De La Frontera'
$ws.Cells.Item(41, 2).Value = 'Amatenango Del Valle'
$ws.Cells.Item(45, 2).Value = 'Benemérito De Las Américas'
$ws.Cells.Item(55, 2).Value = 'Chiapa De Corzo'
$ws.Cells.Item(62, 2).Value = 'Comitán De Domínguez'
$ws.Cells.Item(91, 2).Value = 'Marqués De Comillas'
$ws.Cells.Item(92, 2).Value = 'Mazapa De Madero'
$ws.Cells.Item(96, 2).Value = 'Montecristo De Guerrero'
$ws.Cells.Item(101, 2).Value = 'Ocozocoautla De Espinosa'
$ws.Cells.Item(112, 2).Value = 'Salto De Agua'
$ws.Cells.Item(113, 2).Value = 'San Cristóbal De Las Casas'
$ws.Cells.Item(159, 2).Value = 'Coyame Del Sotol'
$ws.Cells.Item(169, 2).Value = 'Guadalupe Y Calvo'
$ws.Cells.Item(172, 2).Value = 'Hidalgo Del Parral'
$ws.Cells.Item(191, 2).Value = 'San Francisco De Borja'
$ws.Cells.Item(192, 2).Value = 'San Francisco De Conchos'
$ws.Cells.Item(193, 2).Value = 'San Francisco Del Oro'
$ws.Cells.Item(201, 2).Value = 'Valle De Zaragoza'
$ws.Cells.Item(219, 2).Value = 'San Juan De Sabinas'
$ws.Cells.Item(234, 2).Value = 'Villa De Álvarez'
$ws.Cells.Item(236, 1).Value = 'Ciudad De México'
$ws.Cells.Item(240, 2).Value = 'Cuajimalpa De Morelos'
$ws.Cells.Item(255, 2).Value = 'Coneto De Comonfort'
$ws.Cells.Item(269, 2).Value = 'Nombre De Dios'
$ws.Cells.Item(273, 2).Value = 'Pánuco De Coronado'
$ws.Cells.Item(280, 2).Value = 'San Juan De Guadalupe'
$ws.Cells.Item(281, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(282, 2).Value = 'San Luis Del Cordero'
$ws.Cells.Item(283, 2).Value = 'San Pedro Del Gallo'
$ws.Cells.Item(293, 1).Value = 'Estado De México'
$ws.Cells.Item(293, 2).Value = 'Acambay De Ruíz Castañeda'
$ws.Cells.Item(296, 2).Value = 'Almoloya De Alquisiras'
$ws.Cells.Item(297, 2).Value = 'Almoloya De Juárez'
$ws.Cells.Item(298, 2).Value = 'Almoloya Del Río'
$ws.Cells.Item(305, 2).Value = 'Atizapán De Zaragoza'
$ws.Cells.Item(313, 2).Value = 'Chapa De Mota'
$ws.Cells.Item(319, 2).Value = 'Coacalco De Berriozábal'
$ws.Cells.Item(326, 2).Value = 'Ecatepec De Morelos'
$ws.Cells.Item(334, 2).Value = 'Ixtapan De La Sal'
$ws.Cells.Item(335, 2).Value = 'Ixtapan Del Oro'
$ws.Cells.Item(352, 2).Value = 'Naucalpan De Juárez'
$ws.Cells.Item(366, 2).Value = 'San Antonio La Isla'
$ws.Cells.Item(367, 2).Value = 'San Felipe Del Progreso'
$ws.Cells.Item(368, 2).Value = 'San Martín De Las Pirámides'
$ws.Cells.Item(370, 2).Value = 'San Simón De Guerrero'
$ws.Cells.Item(372, 2).Value = 'Soyaniquilpan De Juárez'
$ws.Cells.Item(381, 2).Value = 'Tenango Del Aire'
$ws.Cells.Item(382, 2).Value = 'Tenango Del Valle'
$ws.Cells.Item(396, 2).Value = 'Tlalnepantla De Baz'
$ws.Cells.Item(402, 2).Value = 'Valle De Bravo'
$ws.Cells.Item(403, 2).Value = 'Valle De Chalco Solidaridad'
$ws.Cells.Item(404, 2).Value = 'Villa De Allende'
$ws.Cells.Item(405, 2).Value = 'Villa Del Carbón'
$ws.Cells.Item(419, 2).Value = 'Apaseo El Alto'
$ws.Cells.Item(420, 2).Value = 'Apaseo El Grande'
$ws.Cells.Item(428, 2).Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Cells.Item(432, 2).Value = 'Jaral Del Progreso'
$ws.Cells.Item(440, 2).Value = 'Purísima Del Rincón'
$ws.Cells.Item(444, 2).Value = 'San Diego De La Unión'
$ws.Cells.Item(446, 2).Value = 'San Francisco Del Rincón'
$ws.Cells.Item(448, 2).Value = 'San Luis De La Paz'
$ws.Cells.Item(450, 2).Value = 'Santa Cruz De Juventino Rosas'
$ws.Cells.Item(452, 2).Value = 'Silao De La Victoria'
$ws.Cells.Item(457, 2).Value = 'Valle De Santiago'
$ws.Cells.Item(463, 2).Value = 'Acapulco De Juárez'
$ws.Cells.Item(466, 2).Value = 'Ajuchitlán Del Progreso'
$ws.Cells.Item(467, 2).Value = 'Alcozauca De Guerrero'
$ws.Cells.Item(471, 2).Value = 'Atenango Del Río'
$ws.Cells.Item(472, 2).Value = 'Atlamajalcingo Del Monte'
$ws.Cells.Item(474, 2).Value = 'Atoyac De Álvarez'
$ws.Cells.Item(475, 2).Value = 'Ayutla De Los Libres'
$ws.Cells.Item(478, 2).Value = 'Buenavista De Cuéllar'
$ws.Cells.Item(479, 2).Value = 'Chilapa De Álvarez'
$ws.Cells.Item(480, 2).Value = 'Chilpancingo De Los Bravo'
$ws.Cells.Item(481, 2).Value = 'Coahuayutla De José María Izazaga'
$ws.Cells.Item(486, 2).Value = 'Coyuca De Benítez'
$ws.Cells.Item(487, 2).Value = 'Coyuca De Catalán'
$ws.Cells.Item(491, 2).Value = 'Cuetzala Del Progreso'
$ws.Cells.Item(492, 2).Value = 'Cutzamala De Pinzón'
$ws.Cells.Item(498, 2).Value = 'Huitzuco De Los Figueroa'
$ws.Cells.Item(499, 2).Value = 'Iguala De La Independencia'
$ws.Cells.Item(501, 2).Value = 'Ixcateopan De Cuauhtémoc'
$ws.Cells.Item(502, 2).Value = 'Zihuatanejo De Azueta'
$ws.Cells.Item(504, 2).Value = 'La Unión De Isidoro Montes De Oca'
$ws.Cells.Item(507, 2).Value = 'Mártir De Cuilapan'
$ws.Cells.Item(520, 2).Value = 'Taxco De Alarcón'
$ws.Cells.Item(522, 2).Value = 'Técpan De Galeana'
$ws.Cells.Item(524, 2).Value = 'Tepecoacuilco De Trujano'
$ws.Cells.Item(526, 2).Value = 'Tixtla De Guerrero'
$ws.Cells.Item(530, 2).Value = 'Tlalixtaquilla De Maldonado'
$ws.Cells.Item(531, 2).Value = 'Tlapa De Comonfort'
$ws.Cells.Item(543, 2).Value = 'Agua Blanca De Iturbide'
$ws.Cells.Item(550, 2).Value = 'Atotonilco De Tula'
$ws.Cells.Item(551, 2).Value = 'Atotonilco El Grande'
$ws.Cells.Item(557, 2).Value = 'Cuautepec De Hinojosa'
$ws.Cells.Item(563, 2).Value = 'Huasca De Ocampo'
$ws.Cells.Item(567, 2).Value = 'Huejutla De Reyes'
$ws.Cells.Item(576, 2).Value = 'Mineral De La Reforma'
$ws.Cells.Item(577, 2).Value = 'Mineral Del Chico'
$ws.Cells.Item(578, 2).Value = 'Mineral Del Monte'
$ws.Cells.Item(579, 2).Value = 'Mixquiahuala De Juárez'
$ws.Cells.Item(580, 2).Value = 'Molango De Escamilla'
$ws.Cells.Item(582, 2).Value = 'Nopala De Villagrán'
$ws.Cells.Item(583, 2).Value = 'Omitlán De Juárez'
$ws.Cells.Item(584, 2).Value = 'Pachuca De Soto'
$ws.Cells.Item(587, 2).Value = 'Progreso De Obregón'
$ws.Cells.Item(593, 2).Value = 'Santiago De Anaya'
$ws.Cells.Item(594, 2).Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Cells.Item(598, 2).Value = 'Tenango De Doria'
$ws.Cells.Item(600, 2).Value = 'Tepehuacán De Guerrero'
$ws.Cells.Item(601, 2).Value = 'Tepeji Del Río De Ocampo'
$ws.Cells.Item(604, 2).Value = 'Tezontepec De Aldama'
$ws.Cells.Item(613, 2).Value = 'Tula De Allende'
$ws.Cells.Item(614, 2).Value = 'Tulancingo De Bravo'
$ws.Cells.Item(615, 2).Value = 'Villa De Tezontepec'
$ws.Cells.Item(619, 2).Value = 'Zacualtipán De Ángeles'
$ws.Cells.Item(620, 2).Value = 'Zapotlán De Juárez'
$ws.Cells.Item(625, 2).Value = 'Acatlán De Juárez'
$ws.Cells.Item(626, 2).Value = 'Ahualulco De Mercado'
$ws.Cells.Item(631, 2).Value = 'Atemajac De Brizuela'
$ws.Cells.Item(634, 2).Value = 'Atotonilco El Alto'
$ws.Cells.Item(636, 2).Value = 'Autlán De Navarro'
$ws.Cells.Item(642, 2).Value = 'Cañadas De Obregón'
$ws.Cells.Item(649, 2).Value = 'Concepción De Buenos Aires'
$ws.Cells.Item(650, 2).Value = 'Cuautitlán De García Barragán'
$ws.Cells.Item(659, 2).Value = 'Encarnación De Díaz'
$ws.Cells.Item(666, 2).Value = 'Huejuquilla El Alto'
$ws.Cells.Item(667, 2).Value = 'Ixtlahuacán De Los Membrillos'
$ws.Cells.Item(668, 2).Value = 'Ixtlahuacán Del Río'
$ws.Cells.Item(672, 2).Value = 'Jilotlán De Los Dolores'
$ws.Cells.Item(678, 2).Value = 'La Manzanilla De La Paz'
$ws.Cells.Item(679, 2).Value = 'Lagos De Moreno'
$ws.Cells.Item(687, 2).Value = 'Ojuelos De Jalisco'
$ws.Cells.Item(692, 2).Value = 'San Cristóbal De La Barranca'
$ws.Cells.Item(693, 2).Value = 'San Diego De Alejandría'
$ws.Cells.Item(695, 2).Value = 'San Juan De Los Lagos'
$ws.Cells.Item(696, 2).Value = 'San Juanito De Escobedo'
$ws.Cells.Item(699, 2).Value = 'San Martín De Bolaños'
$ws.Cells.Item(701, 2).Value = 'San Miguel El Alto'
$ws.Cells.Item(702, 2).Value = 'San Sebastián Del Oeste'
$ws.Cells.Item(703, 2).Value = 'Santa María De Los Ángeles'
$ws.Cells.Item(704, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(707, 2).Value = 'Talpa De Allende'
$ws.Cells.Item(708, 2).Value = 'Tamazula De Gordiano'
$ws.Cells.Item(711, 2).Value = 'Techaluta De Montenegro'
$ws.Cells.Item(715, 2).Value = 'Teocuitatlán De Corona'
$ws.Cells.Item(716, 2).Value = 'Tepatitlán De Morelos'
$ws.Cells.Item(719, 2).Value = 'Tizapán El Alto'
$ws.Cells.Item(720, 2).Value = 'Tlajomulco De Zúñiga'
$ws.Cells.Item(732, 2).Value = 'Unión De San Antonio'
$ws.Cells.Item(733, 2).Value = 'Unión De Tula'
$ws.Cells.Item(734, 2).Value = 'Valle De Guadalupe'
$ws.Cells.Item(735, 2).Value = 'Valle De Juárez'
$ws.Cells.Item(740, 2).Value = 'Yahualica De González Gallo'
$ws.Cells.Item(741, 2).Value = 'Zacoalco De Torres'
$ws.Cells.Item(744, 2).Value = 'Zapotitlán De Vadillo'
$ws.Cells.Item(745, 2).Value = 'Zapotlán Del Rey'
$ws.Cells.Item(746, 2).Value = 'Zapotlán El Grande'
$ws.Cells.Item(772, 2).Value = 'Coalcomán De Vázquez Pallares'
$ws.Cells.Item(774, 2).Value = 'Cojumatlán De Régules'
$ws.Cells.Item(841, 2).Value = 'Tiquicheo De Nicolás Romero'
$ws.Cells.Item(867, 2).Value = 'Coatlán Del Río'
$ws.Cells.Item(879, 2).Value = 'Puente De Ixtla'
$ws.Cells.Item(885, 2).Value = 'Tetela Del Volcán'
$ws.Cells.Item(887, 2).Value = 'Tlaltizapán De Zapata'
$ws.Cells.Item(895, 2).Value = 'Zacualpan De Amilpas'
$ws.Cells.Item(899, 2).Value = 'Amatlán De Cañas'
$ws.Cells.Item(900, 2).Value = 'Bahía De Banderas'
$ws.Cells.Item(904, 2).Value = 'Ixtlán Del Río'
$ws.Cells.Item(911, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(940, 2).Value = 'Mier Y Noriega'
$ws.Cells.Item(944, 2).Value = 'San Nicolás De Los Garza'
$ws.Cells.Item(951, 2).Value = 'Acatlán De Pérez Figueroa'
$ws.Cells.Item(958, 2).Value = 'Ayoquezco De Aldama'
$ws.Cells.Item(962, 2).Value = 'Capulálpam De Méndez'
$ws.Cells.Item(964, 2).Value = 'Chalcatongo De Hidalgo'
$ws.Cells.Item(965, 2).Value = 'Chiquihuitlán De Benito Juárez'
$ws.Cells.Item(968, 2).Value = 'Coicoyán De Las Flores'
$ws.Cells.Item(971, 2).Value = 'Constancia Del Rosario'
$ws.Cells.Item(974, 2).Value = 'Cuilápam De Guerrero'
$ws.Cells.Item(975, 2).Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Cells.Item(976, 2).Value = 'El Barrio De La Soledad'
$ws.Cells.Item(978, 2).Value = 'Eloxochitlán De Flores Magón'
$ws.Cells.Item(979, 2).Value = 'Fresnillo De Trujano'
$ws.Cells.Item(980, 2).Value = 'Guadalupe De Ramírez'
$ws.Cells.Item(982, 2).Value = 'Guelatao De Juárez'
$ws.Cells.Item(983, 2).Value = 'Guevea De Humboldt'
$ws.Cells.Item(984, 2).Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Cells.Item(985, 2).Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Cells.Item(986, 2).Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Cells.Item(987, 2).Value = 'Huautla De Jiménez'
$ws.Cells.Item(989, 2).Value = 'Ixtlán De Juárez'
$ws.Cells.Item(990, 2).Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Cells.Item(1004, 2).Value = 'Magdalena Yodocono De Porfirio Díaz'
$ws.Cells.Item(1006, 2).Value = 'Mariscala De Juárez'
$ws.Cells.Item(1007, 2).Value = 'Mártires De Tacubaya'
$ws.Cells.Item(1009, 2).Value = 'Mazatlán Villa De Flores'
$ws.Cells.Item(1011, 2).Value = 'Miahuatlán De Porfirio Díaz'
$ws.Cells.Item(1012, 2).Value = 'Mixistlán De La Reforma'
$ws.Cells.Item(1016, 2).Value = 'Nejapa De Madero'
$ws.Cells.Item(1018, 2).Value = 'Oaxaca De Juárez'
$ws.Cells.Item(1019, 2).Value = 'Ocotlán De Morelos'
$ws.Cells.Item(1020, 2).Value = 'Pinotepa De Don Luis'
$ws.Cells.Item(1022, 2).Value = 'Putla Villa De Guerrero'
$ws.Cells.Item(1023, 2).Value = 'Reforma De Pineda'
$ws.Cells.Item(1025, 2).Value = 'Rojas De Cuauhtémoc'
$ws.Cells.Item(1030, 2).Value = 'San Agustín De Las Juntas'
$ws.Cells.Item(1051, 2).Value = 'San Antonino El Alto'
$ws.Cells.Item(1053, 2).Value = 'San Antonio De La Cal'
$ws.Cells.Item(1060, 2).Value = 'San Baltazar Yatzachi El Bajo'
$ws.Cells.Item(1076, 2).Value = 'San Dionisio Del Mar'
$ws.Cells.Item(1080, 2).Value = 'San Felipe Jalapa De Díaz'
$ws.Cells.Item(1087, 2).Value = 'San Francisco Del Mar'
$ws.Cells.Item(1113, 2).Value = 'San José Del Peñasco'
$ws.Cells.Item(1114, 2).Value = 'San José Del Progreso'
$ws.Cells.Item(1125, 2).Value = 'San Juan Bautista Lo De Soto'
$ws.Cells.Item(1139, 2).Value = 'San Juan De Los Cués'
$ws.Cells.Item(1140, 2).Value = 'San Juan Del Estado'
$ws.Cells.Item(1141, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(1183, 2).Value = 'San Martín De Los Cansecos'
$ws.Cells.Item(1192, 2).Value = 'San Mateo Del Mar'
$ws.Cells.Item(1209, 2).Value = 'San Miguel Del Puerto'
$ws.Cells.Item(1210, 2).Value = 'San Miguel Del Río'
$ws.Cells.Item(1212, 2).Value = 'San Miguel El Grande'
$ws.Cells.Item(1235, 2).Value = 'San Pablo Villa De Mitla'
$ws.Cells.Item(1241, 2).Value = 'San Pedro El Alto'
$ws.Cells.Item(1266, 2).Value = 'San Pedro Y San Pablo Ayutla'
$ws.Cells.Item(1267, 2).Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Cells.Item(1268, 2).Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Cells.Item(1286, 2).Value = 'Santa Ana Del Valle'
$ws.Cells.Item(1304, 2).Value = 'Santa Cruz De Bravo'
$ws.Cells.Item(1309, 2).Value = 'Santa Cruz Tacache De Mina'
$ws.Cells.Item(1316, 2).Value = 'Santa Inés De Zaragoza'
$ws.Cells.Item(1317, 2).Value = 'Santa Inés Del Monte'
$ws.Cells.Item(1319, 2).Value = 'Santa Lucía Del Camino'
$ws.Cells.Item(1333, 2).Value = 'Santa María Del Tule'
$ws.Cells.Item(1341, 2).Value = 'Santa María Jalapa Del Marqués'
$ws.Cells.Item(1343, 2).Value = 'Santa María La Asunción'
$ws.Cells.Item(1383, 2).Value = 'Santiago Del Río'
$ws.Cells.Item(1426, 2).Value = 'Santo Domingo De Morelos'
$ws.Cells.Item(1450, 2).Value = 'Sitio De Xitlapehua'
$ws.Cells.Item(1452, 2).Value = 'Tamazulápam Del Espíritu Santo'
$ws.Cells.Item(1453, 2).Value = 'Tanetze De Zaragoza'
$ws.Cells.Item(1455, 2).Value = 'Tataltepec De Valdés'
$ws.Cells.Item(1456, 2).Value = 'Teococuilco De Marcos Pérez'
$ws.Cells.Item(1457, 2).Value = 'Teotitlán De Flores Magón'
$ws.Cells.Item(1458, 2).Value = 'Teotitlán Del Valle'
$ws.Cells.Item(1460, 2).Value = 'Tepelmeme Villa De Morelos'
$ws.Cells.Item(1461, 2).Value = 'Tezoatlán De Segura Y Luna'
$ws.Cells.Item(1462, 2).Value = 'Tlacolula De Matamoros'
$ws.Cells.Item(1464, 2).Value = 'Tlalixtac De Cabrera'
$ws.Cells.Item(1465, 2).Value = 'Totontepec Villa De Morelos'
$ws.Cells.Item(1469, 2).Value = 'Villa De Chilapa De Díaz'
$ws.Cells.Item(1470, 2).Value = 'Villa De Etla'
$ws.Cells.Item(1471, 2).Value = 'Villa De Tamazulápam Del Progreso'
$ws.Cells.Item(1472, 2).Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Cells.Item(1473, 2).Value = 'Villa De Zaachila'
$ws.Cells.Item(1476, 2).Value = 'Villa Sola De Vega'
$ws.Cells.Item(1477, 2).Value = 'Villa Talea De Castro'
$ws.Cells.Item(1478, 2).Value = 'Villa Tejúpam De La Unión'
$ws.Cells.Item(1481, 2).Value = 'Yutanduchi De Guerrero'
$ws.Cells.Item(1482, 2).Value = 'Zapotitlán Del Río'
$ws.Cells.Item(1485, 2).Value = 'Zimatlán De Álvarez'
$ws.Cells.Item(1512, 2).Value = 'Ayotoxco De Guerrero'
$ws.Cells.Item(1518, 2).Value = 'Chalchicomula De Sesma'
$ws.Cells.Item(1528, 2).Value = 'Chila De La Sal'
$ws.Cells.Item(1539, 2).Value = 'Cuapiaxtla De Madero'
$ws.Cells.Item(1543, 2).Value = 'Cuayuca De Andrade'
$ws.Cells.Item(1544, 2).Value = 'Cuetzalan Del Progreso'
$ws.Cells.Item(1560, 2).Value = 'Huehuetlán El Chico'
$ws.Cells.Item(1561, 2).Value = 'Huehuetlán El Grande'
$ws.Cells.Item(1566, 2).Value = 'Huitzilan De Serdán'
$ws.Cells.Item(1568, 2).Value = 'Ixcamilpa De Guerrero'
$ws.Cells.Item(1572, 2).Value = 'Izúcar De Matamoros'
$ws.Cells.Item(1582, 2).Value = 'Los Reyes De Juárez'
$ws.Cells.Item(1583, 2).Value = 'Mazapiltepec De Juárez'
$ws.Cells.Item(1596, 2).Value = 'Palmar De Bravo'
$ws.Cells.Item(1606, 2).Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Cells.Item(1623, 2).Value = 'San Nicolás De Los Ranchos'
$ws.Cells.Item(1627, 2).Value = 'San Salvador El Seco'
$ws.Cells.Item(1628, 2).Value = 'San Salvador El Verde'
$ws.Cells.Item(1637, 2).Value = 'Tecali De Herrera'
$ws.Cells.Item(1645, 2).Value = 'Tepanco De López'
$ws.Cells.Item(1646, 2).Value = 'Tepango De Rodríguez'
$ws.Cells.Item(1647, 2).Value = 'Tepatlaxco De Hidalgo'
$ws.Cells.Item(1652, 2).Value = 'Tepexi De Rodríguez'
$ws.Cells.Item(1654, 2).Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Cells.Item(1655, 2).Value = 'Tetela De Ocampo'
$ws.Cells.Item(1656, 2).Value = 'Teteles De Avila Castillo'
$ws.Cells.Item(1661, 2).Value = 'Tlacotepec De Benito Juárez'
$ws.Cells.Item(1673, 2).Value = 'Totoltepec De Guerrero'
$ws.Cells.Item(1675, 2).Value = 'Tuzamapan De Galeana'
$ws.Cells.Item(1679, 2).Value = 'Xayacatlán De Bravo'
$ws.Cells.Item(1685, 2).Value = 'Xochitlán De Vicente Suárez'
$ws.Cells.Item(1700, 2).Value = 'Amealco De Bonfil'
$ws.Cells.Item(1702, 2).Value = 'Cadereyta De Montes'
$ws.Cells.Item(1708, 2).Value = 'Jalpan De Serra'
$ws.Cells.Item(1709, 2).Value = 'Landa De Matamoros'
$ws.Cells.Item(1712, 2).Value = 'Pinal De Amoles'
$ws.Cells.Item(1715, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(1729, 2).Value = 'Armadillo De Los Infante'
$ws.Cells.Item(1730, 2).Value = 'Axtla De Terrazas'
$ws.Cells.Item(1735, 2).Value = 'Ciudad Del Maíz'
$ws.Cells.Item(1745, 2).Value = 'Mexquitic De Carmona'
$ws.Cells.Item(1751, 2).Value = 'San Ciro De Acosta'
$ws.Cells.Item(1756, 2).Value = 'Santa María Del Río'
$ws.Cells.Item(1758, 2).Value = 'Soledad De Graciano Sánchez'
$ws.Cells.Item(1765, 2).Value = 'Tanquián De Escobedo'
$ws.Cells.Item(1768, 2).Value = 'Villa De Arista'
$ws.Cells.Item(1769, 2).Value = 'Villa De Arriaga'
$ws.Cells.Item(1770, 2).Value = 'Villa De Guadalupe'
$ws.Cells.Item(1771, 2).Value = 'Villa De Ramos'
$ws.Cells.Item(1772, 2).Value = 'Villa De Reyes'
$ws.Cells.Item(1832, 2).Value = 'Nacozari De García'
$ws.Cells.Item(1841, 2).Value = 'San Felipe De Jesús'
$ws.Cells.Item(1844, 2).Value = 'San Miguel De Horcasitas'
$ws.Cells.Item(1863, 2).Value = 'Jalpa De Méndez'
$ws.Cells.Item(1894, 2).Value = 'Soto La Marina'
$ws.Cells.Item(1901, 2).Value = 'Acuamanala De Miguel Hidalgo'
$ws.Cells.Item(1903, 2).Value = 'Amaxac De Guerrero'
$ws.Cells.Item(1904, 2).Value = 'Apetatitlán De Antonio Carvajal'
$ws.Cells.Item(1910, 2).Value = 'Contla De Juan Cuamatzi'
$ws.Cells.Item(1917, 2).Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Cells.Item(1920, 2).Value = 'Muñoz De Domingo Arenas'
$ws.Cells.Item(1921, 2).Value = 'Nanacamilpa De Mariano Arista'
$ws.Cells.Item(1924, 2).Value = 'Papalotla De Xicohténcatl'
$ws.Cells.Item(1927, 2).Value = 'San Pablo Del Monte'
$ws.Cells.Item(1934, 2).Value = 'Tepetitla De Lardizábal'
$ws.Cells.Item(1937, 2).Value = 'Tetla De La Solidaridad'
$ws.Cells.Item(1949, 2).Value = 'Ziltlaltépec De Trinidad Sánchez Santos'
$ws.Cells.Item(1959, 2).Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Cells.Item(1963, 2).Value = 'Amatlán De Los Reyes'
$ws.Cells.Item(1975, 2).Value = 'Boca Del Río'
$ws.Cells.Item(1980, 2).Value = 'Castillo De Teayo'
$ws.Cells.Item(1990, 2).Value = 'Chinampa De Gorostiza'
$ws.Cells.Item(2001, 2).Value = 'Cosamaloapan De Carpio'
$ws.Cells.Item(2002, 2).Value = 'Cosautlán De Carvajal'
$ws.Cells.Item(2019, 2).Value = 'Hueyapan De Ocampo'
$ws.Cells.Item(2020, 2).Value = 'Huiloapan De Cuauhtémoc'
$ws.Cells.Item(2021, 2).Value = 'Ignacio De La Llave'
$ws.Cells.Item(2025, 2).Value = 'Ixhuacán De Los Reyes'
$ws.Cells.Item(2026, 2).Value = 'Ixhuatlán De Madero'
$ws.Cells.Item(2027, 2).Value = 'Ixhuatlán Del Café'
$ws.Cells.Item(2028, 2).Value = 'Ixhuatlán Del Sureste'
$ws.Cells.Item(2039, 2).Value = 'Juchique De Ferrer'
$ws.Cells.Item(2044, 2).Value = 'Las Vigas De Ramírez'
$ws.Cells.Item(2045, 2).Value = 'Lerdo De Tejada'
$ws.Cells.Item(2050, 2).Value = 'Martínez De La Torre'
$ws.Cells.Item(2052, 2).Value = 'Medellín De Bravo'
$ws.Cells.Item(2055, 2).Value = 'Mixtla De Altamirano'
$ws.Cells.Item(2057, 2).Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Cells.Item(2068, 2).Value = 'Ozuluama De Mascareñas'
$ws.Cells.Item(2072, 2).Value = 'Paso De Ovejas'
$ws.Cells.Item(2073, 2).Value = 'Paso Del Macho'
$ws.Cells.Item(2076, 2).Value = 'Poza Rica De Hidalgo'
$ws.Cells.Item(2086, 2).Value = 'Sayula De Alemán'
$ws.Cells.Item(2090, 2).Value = 'Soledad De Doblado'
$ws.Cells.Item(2097, 2).Value = 'Tatahuicapan De Juárez'
$ws.Cells.Item(2118, 2).Value = 'Tlacotepec De Mejía'
$ws.Cells.Item(2132, 2).Value = 'Vega De Alatorre'
$ws.Cells.Item(2143, 2).Value = 'Zontecomatlán De López Y Fuentes'
$ws.Cells.Item(2144, 2).Value = 'Zozocolco De Hidalgo'
$ws.Cells.Item(2221, 2).Value = 'Cañitas De Felipe Pescador'
$ws.Cells.Item(2223, 2).Value = 'Concepción Del Oro'
$ws.Cells.Item(2225, 2).Value = 'El Plateado De Joaquín Amaro'
$ws.Cells.Item(2235, 2).Value = 'Jiménez Del Teul'
$ws.Cells.Item(2241, 2).Value = 'Mezquital Del Oro'
$ws.Cells.Item(2246, 2).Value = 'Moyahua De Estrada'
$ws.Cells.Item(2247, 2).Value = 'Nochistlán De Mejía'
$ws.Cells.Item(2248, 2).Value = 'Noria De Ángeles'
$ws.Cells.Item(2259, 2).Value = 'Teúl De González Ortega'
$ws.Cells.Item(2260, 2).Value = 'Tlaltenango De Sánchez Román'
$ws.Cells.Item(2262, 2).Value = 'Trinidad García De La Cadena'
$ws.Cells.Item(2265, 2).Value = 'Villa De Cos'

# --- Special-case capitalization fix ---
$ws.Cells.Item(941, 2).Value = 'Montemorelos'

# --- Floating point recompute touch-ups (percentage column) ---
$ws.Cells.Item(156, 4).Value = 0.0009047158312704976
$ws.Cells.Item(741, 4).Value = 0.0009235640777552996
$ws.Cells.Item(1258, 4).Value = 0.000917281328927032
$ws.Cells.Item(1489, 4).Value = 0.000917281328927032
$ws.Cells.Item(1493, 4).Value = 0.0009047158312704976
$ws.Cells.Item(1517, 4).Value = 0.0009486950730683688
$ws.Cells.Item(1525, 4).Value = 0.000917281328927032
$ws.Cells.Item(1786, 4).Value = 0.000917281328927032
$ws.Cells.Item(1787, 4).Value = 0.000917281328927032

# --- Drop the trailing footnote/metadata rows (2274:2278) so the used range shrinks to D2272 ---
$ws.Range("A2274:D2278").ClearContents()
